$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the claim/annulment numbers (rows 2-4); these are text values
# (leading-zero numeric strings), so use an apostrophe prefix to force
# text interpretation, matching the original quotePrefix-styled cells.
$ws.Range("B2").Value = "'0420194406717"
$ws.Range("C2").Value = "'4500996"

$ws.Range("B3").Value = "'1120194100412"
$ws.Range("C3").Value = "'4500205"

$ws.Range("B4").Value = "'1220194200667"
$ws.Range("C4").Value = "'4500278"

# Update the active selection to B5
$ws.Range("B5").Select()
